# Update "Time" observed-at derived data: recomputed color-sampling
# statistics (K, L-N, O, P-S, T, U-X, Y, Z columns) for rows 2-15 on Sheet1
# to reflect the refreshed "time observed at" analysis run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = '[''#2a4b13'', ''#bccbdd'', ''#a3b2bd'']'
$ws.Range("L2").Value = 188.3416957360965
$ws.Range("M2").Value = 202.8464710265135
$ws.Range("N2").Value = 220.915763538457
$ws.Range("O2").Value = 'darkgreen'
$ws.Range("Q2").Value = 42.29902360588153
$ws.Range("R2").Value = 75.41250902988652
$ws.Range("S2").Value = 19.38843707034121
$ws.Range("T2").Value = 'lightsteelblue'
$ws.Range("V2").Value = 163.2768482768487
$ws.Range("W2").Value = 177.9895133228466
$ws.Range("X2").Value = 189.2472472472464
$ws.Range("Y2").Value = 'darkgray'
$ws.Range("K3").Value = '[''#b7c6db'', ''#5c6c44'', ''#393117'']'
$ws.Range("O3").Value = 'lightsteelblue'
$ws.Range("P3").Value = 0.6691058201058201
$ws.Range("Q3").Value = 57.35710091176901
$ws.Range("R3").Value = 48.51392365940644
$ws.Range("S3").Value = 22.7506722299506
$ws.Range("T3").Value = 'darkolivegreen'
$ws.Range("U3").Value = 0.171973544973545
$ws.Range("V3").Value = 92.11828455904576
$ws.Range("W3").Value = 107.7036754714501
$ws.Range("X3").Value = 67.50449076990198
$ws.Range("Y3").Value = 'darkslategray'
$ws.Range("Z3").Value = 0.1589206349206349
$ws.Range("K4").Value = '[''#adbcd2'', ''#233827'', ''#9aa8bd'']'
$ws.Range("L4").Value = 172.7732299764028
$ws.Range("M4").Value = 187.6358839623083
$ws.Range("N4").Value = 210.2517768995839
$ws.Range("O4").Value = 'lightsteelblue'
$ws.Range("Q4").Value = 153.8308242811499
$ws.Range("R4").Value = 167.9830031948883
$ws.Range("S4").Value = 188.8616996805116
$ws.Range("T4").Value = 'darkslategray'
$ws.Range("V4").Value = 35.32824211570215
$ws.Range("W4").Value = 56.05699335311498
$ws.Range("X4").Value = 39.02866166969977
$ws.Range("Y4").Value = 'darkgray'
$ws.Range("K5").Value = '[''#acbad2'', ''#3f3e2a'', ''#97a4ba'']'
$ws.Range("L5").Value = 171.6506978134842
$ws.Range("M5").Value = 186.0438363827481
$ws.Range("N5").Value = 209.8670441064442
$ws.Range("O5").Value = 'lightsteelblue'
$ws.Range("Q5").Value = 151.0836809013941
$ws.Range("R5").Value = 163.9464797061756
$ws.Range("S5").Value = 185.5989946464142
$ws.Range("T5").Value = 'darkslategray'
$ws.Range("V5").Value = 62.50738639268461
$ws.Range("W5").Value = 61.79379440779776
$ws.Range("X5").Value = 41.63170270407647
$ws.Range("Y5").Value = 'darkgray'
$ws.Range("O6").Value = 'lightslategray'
$ws.Range("S6").Value = 197.2461172689395
$ws.Range("T6").Value = 'darkgray'
$ws.Range("V6").Value = 52.66094069530308
$ws.Range("Y6").Value = 'darkslategray'
$ws.Range("K7").Value = '[''#828fa9'', ''#384935'', ''#8e9bbc'']'
$ws.Range("O7").Value = 'lightslategray'
$ws.Range("T7").Value = 'darkslategray'
$ws.Range("X7").Value = 52.70864395292821
$ws.Range("Y7").Value = 'darkgray'
$ws.Range("K8").Value = '[''#547528'', ''#dedad6'', ''#365b10'']'
$ws.Range("L8").Value = 221.8846870062522
$ws.Range("M8").Value = 218.1225171168326
$ws.Range("N8").Value = 213.5128846587863
$ws.Range("O8").Value = 'darkolivegreen'
$ws.Range("P8").Value = 0.5656010638297873
$ws.Range("Q8").Value = 53.75551149508335
$ws.Range("R8").Value = 90.61639050180048
$ws.Range("S8").Value = 15.63958520195835
$ws.Range("T8").Value = 'gainsboro'
$ws.Range("U8").Value = 0.3361010638297872
$ws.Range("V8").Value = 84.49048170832825
$ws.Range("W8").Value = 117.3252220934754
$ws.Range("X8").Value = 39.85057661535133
$ws.Range("Y8").Value = 'darkolivegreen'
$ws.Range("Z8").Value = 0.09829787234042553
$ws.Range("K9").Value = '[''#244a0f'', ''#dbdad9'', ''#849578'']'
$ws.Range("L9").Value = 219.1860148440787
$ws.Range("M9").Value = 217.7420916382349
$ws.Range("N9").Value = 217.2547061112169
$ws.Range("O9").Value = 'darkgreen'
$ws.Range("P9").Value = 0.7033578947368421
$ws.Range("Q9").Value = 36.43105277114961
$ws.Range("R9").Value = 73.84079630260167
$ws.Range("S9").Value = 14.80382240404404
$ws.Range("T9").Value = 'gainsboro'
$ws.Range("U9").Value = 0.2779105263157895
$ws.Range("V9").Value = 132.2873239436621
$ws.Range("W9").Value = 149.1107042253521
$ws.Range("X9").Value = 120.155492957746
$ws.Range("Z9").Value = 0.01873157894736842
$ws.Range("K10").Value = '[''#d3d1ce'', ''#1b4009'', ''#4d6a28'']'
$ws.Range("O10").Value = 'lightgray'
$ws.Range("R10").Value = 64.45220293724373
$ws.Range("T10").Value = 'darkgreen'
$ws.Range("Y10").Value = 'darkolivegreen'
$ws.Range("K11").Value = '[''#d7d2cd'', ''#54643a'', ''#3b4c1b'']'
$ws.Range("L11").Value = 214.5962834987609
$ws.Range("M11").Value = 210.341278400136
$ws.Range("N11").Value = 204.9647423905637
$ws.Range("O11").Value = 'lightgray'
$ws.Range("Q11").Value = 59.33412069998509
$ws.Range("R11").Value = 76.39573491250985
$ws.Range("S11").Value = 27.18354514606446
$ws.Range("T11").Value = 'darkolivegreen'
$ws.Range("U11").Value = 0.1653079019073569
$ws.Range("V11").Value = 84.44690911213885
$ws.Range("W11").Value = 100.3524604518122
$ws.Range("X11").Value = 57.54270036776209
$ws.Range("Y11").Value = 'darkolivegreen'
$ws.Range("Z11").Value = 0.09186376021798365
$ws.Range("L12").Value = 201.3772533818934
$ws.Range("M12").Value = 199.7616357775198
$ws.Range("N12").Value = 197.0023373815555
$ws.Range("O12").Value = 'silver'
$ws.Range("Q12").Value = 183.305439533167
$ws.Range("R12").Value = 182.3350403757912
$ws.Range("S12").Value = 180.585140255935
$ws.Range("T12").Value = 'darkolivegreen'
$ws.Range("Y12").Value = 'silver'
$ws.Range("L13").Value = 209.4411456859088
$ws.Range("M13").Value = 208.0204859488259
$ws.Range("N13").Value = 208.4494137937803
$ws.Range("O13").Value = 'lightgray'
$ws.Range("P13").Value = 0.8300765027322404
$ws.Range("S13").Value = 19.91608179370846
$ws.Range("T13").Value = 'darkgreen'
$ws.Range("V13").Value = 159.8989856297552
$ws.Range("W13").Value = 164.1223584108203
$ws.Range("X13").Value = 157.4387151310225
$ws.Range("Y13").Value = 'darkgray'
$ws.Range("Z13").Value = 0.02588524590163934
$ws.Range("K14").Value = '[''#d1cdc8'', ''#22410c'', ''#7f8c78'']'
$ws.Range("L14").Value = 209.0773185863931
$ws.Range("M14").Value = 205.439694473329
$ws.Range("N14").Value = 200.4629688874748
$ws.Range("O14").Value = 'lightgray'
$ws.Range("P14").Value = 0.9164659685863874
$ws.Range("Q14").Value = 33.56708707901527
$ws.Range("R14").Value = 65.45679592111628
$ws.Range("S14").Value = 12.02287669394016
$ws.Range("T14").Value = 'darkgreen'
$ws.Range("U14").Value = 0.07803141361256545
$ws.Range("V14").Value = 127.1730038022816
$ws.Range("W14").Value = 139.7186311787071
$ws.Range("X14").Value = 119.6492395437262
$ws.Range("K15").Value = '[''#27470d'', ''#d6d2d0'', ''#4e6b25'']'
$ws.Range("L15").Value = 214.1110750982702
$ws.Range("M15").Value = 210.4964347286406
$ws.Range("N15").Value = 208.0095441693535
$ws.Range("O15").Value = 'darkgreen'
$ws.Range("P15").Value = 0.7513264248704663
$ws.Range("Q15").Value = 39.20651412561173
$ws.Range("R15").Value = 70.78006070510541
$ws.Range("S15").Value = 12.66495447114249
$ws.Range("T15").Value = 'lightgray'
$ws.Range("U15").Value = 0.1779533678756477
$ws.Range("V15").Value = 78.13164796736125
$ws.Range("W15").Value = 107.0761328864903
$ws.Range("X15").Value = 36.76540871339429
$ws.Range("Y15").Value = 'darkolivegreen'
$ws.Range("Z15").Value = 0.07072020725388602
